$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Location 1: title paragraph
#   "AR ${AR}: Install ..."  ->  "Recommendation ${REC}: Install ..."
# We must be careful to:
#   * only retarget the bare word "AR", not the "AR" inside the "${AR}" field token
#   * leave the bookmark (_Toc286328876), which sits between "AR "/"Recommendation "
#     and "${AR}"/"${REC}", exactly where it is
# so we use small, explicit character-range edits (paragraph-relative offsets)
# instead of a document-wide Find/Replace.
# ------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    $pText = $p.Range.Text
    $wordIdx = $pText.IndexOf("AR `${AR}")
    if ($wordIdx -ge 0) {
        $pStart = $p.Range.Start

        # 1a. Replace the standalone word "AR" (not the "AR" inside "${AR}") with "Recommendation".
        $absStart = $pStart + $wordIdx
        $rWord = $d.Range($absStart, $absStart + 2)
        $rWord.Text = "Recommendation"

        # 1b. Re-locate "${AR}" (its position shifted after the insert above) and turn it into "${REC}".
        $p2 = $d.Paragraphs.Item($i)
        $p2Start = $p2.Range.Start
        $p2Text = $p2.Range.Text
        $tokenIdx = $p2Text.IndexOf("`${AR}")
        if ($tokenIdx -ge 0) {
            $absStart2 = $p2Start + $tokenIdx
            $rToken = $d.Range($absStart2, $absStart2 + 5)
            $rToken.Text = "`${REC}"
        }
        break
    }
}

# ------------------------------------------------------------------
# Location 2: body paragraph
#   "...annual electricity savings for this AR is ${ES}..."
#     -> "...annual electricity savings for this recommendation is ${ES}..."
# ------------------------------------------------------------------
$paraCount = $d.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $p = $d.Paragraphs.Item($i)
    $pText = $p.Range.Text
    $localIdx = $pText.IndexOf("for this AR is")
    if ($localIdx -ge 0) {
        $pStart = $p.Range.Start
        $arLocalIdx = $localIdx + "for this ".Length
        $absIdx = $pStart + $arLocalIdx
        $rAR = $d.Range($absIdx, $absIdx + 2)
        $rAR.Text = "recommendation"
        break
    }
}
